$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits ---
$ws.Range("D2").Value = 44454
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 31000
$ws.Range("P2").Value = 30500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 3050
$ws.Range("T2").Value = 10

# --- Row 3 edits ---
$ws.Range("D3").Value = 44446
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 3200
$ws.Range("O3").Value = 3300
$ws.Range("P3").Value = 3250
$ws.Range("Q3").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R3").Value = 'Provincia del Elquí'
$ws.Range("S3").Value = 3250
$ws.Range("T3").Value = 1

# --- Row 4 edits ---
$ws.Range("D4").Value = 44160
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = '$/bandeja 8 kilos'
$ws.Range("R4").Value = 'Provincia de Limarí'
$ws.Range("S4").Value = 2188
$ws.Range("T4").Value = 8

# --- New row 5 (weekly entry added below; mirrors the price data that used
#     to live in row 3 before that row was overwritten with fresher figures) ---
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 'Vega Monumental Concepción'
$ws.Range("C5").Value = 'Bíobío'
$ws.Range("D5").Value = 44160
$ws.Range("D5").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = 'Otros'
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = 'Chirimoya'
$ws.Range("K5").Value = 'Cultivar IV Región'
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = '$/bandeja 8 kilos'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 1875
$ws.Range("T5").Value = 8
